$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("AB2:AK5").ClearContents()
$ws.Range("AM2:AM5").ClearContents()
